# Update "想去人数" (number of people wanting to attend) counts
# across all four sheets, per gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 116
$ws.Range("F3").Value = 1278
$ws.Range("F4").Value = 915
$ws.Range("F5").Value = 953
$ws.Range("F6").Value = 1706
$ws.Range("F7").Value = 371
$ws.Range("F8").Value = 1140
$ws.Range("F10").Value = 4
$ws.Range("F13").Value = 29
$ws.Range("F15").Value = 630
$ws.Range("F16").Value = 127
$ws.Range("F17").Value = 82
$ws.Range("F18").Value = 22
$ws.Range("F21").Value = 91
$ws.Range("F22").Value = 641
$ws.Range("F23").Value = 12
$ws.Range("F24").Value = 624
$ws.Range("F25").Value = 126
$ws.Range("F26").Value = 29
$ws.Range("F28").Value = 294
$ws.Range("F29").Value = 93
$ws.Range("F30").Value = 23
$ws.Range("F31").Value = 244

# Sheet: 演出 (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 306
$ws.Range("F7").Value = 236
$ws.Range("F11").Value = 109
$ws.Range("F12").Value = 20

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 296

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 296
$ws.Range("F3").Value = 116
$ws.Range("F4").Value = 1278
$ws.Range("F5").Value = 915
$ws.Range("F6").Value = 953
$ws.Range("F7").Value = 1706
$ws.Range("F8").Value = 371
$ws.Range("F9").Value = 1140
$ws.Range("F12").Value = 4
$ws.Range("F15").Value = 29
$ws.Range("F17").Value = 630
$ws.Range("F18").Value = 127
$ws.Range("F19").Value = 82
$ws.Range("F21").Value = 22
$ws.Range("F22").Value = 306
$ws.Range("F27").Value = 236
$ws.Range("F28").Value = 236
$ws.Range("F29").Value = 91
$ws.Range("F30").Value = 641
$ws.Range("F31").Value = 12
$ws.Range("F32").Value = 624
$ws.Range("F33").Value = 126
$ws.Range("F34").Value = 29
$ws.Range("F36").Value = 294
$ws.Range("F39").Value = 93
$ws.Range("F40").Value = 23
$ws.Range("F41").Value = 244
$ws.Range("F43").Value = 109
$ws.Range("F44").Value = 109
$ws.Range("F47").Value = 20
